$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.317.18'
$ws.Range("E2").Value = '  +3.67%  '

$ws.Range("D3").Value = '3.117.08'
$ws.Range("E3").Value = '  +1.63%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '219.09'
$ws.Range("E5").Value = '  +3.84%  '

$ws.Range("D6").Value = '622.99'
$ws.Range("E6").Value = '  +0.74%  '

$ws.Range("D7").Value = '0.992'
$ws.Range("E7").Value = '  +24.36%  '

$ws.Range("D8").Value = '0.377'
$ws.Range("E8").Value = '  +1.93%  '

$ws.Range("D10").Value = '3.113.22'
$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("D11").Value = '0.718'
$ws.Range("E11").Value = '  +20.81%  '

$ws.Range("E12").Value = '  +6.01%  '

$ws.Range("E13").Value = '  +7.34%  '

$ws.Range("D14").Value = '34.69'
$ws.Range("E14").Value = '  +8.46%  '

$ws.Range("E15").Value = '  +2.68%  '

$ws.Range("D16").Value = '91.025.44'
$ws.Range("E16").Value = '  +3.65%  '

$ws.Range("D17").Value = '3.691.39'
$ws.Range("E17").Value = '  +1.70%  '

$ws.Range("D18").Value = '3.118.03'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("D19").Value = '3.78'
$ws.Range("E19").Value = '  +13.39%  '

$ws.Range("E20").Value = '  +6.94%  '

$ws.Range("D21").Value = '14.11'
$ws.Range("E21").Value = '  +5.71%  '

$ws.Range("D22").Value = '436.19'
$ws.Range("E22").Value = '  +3.79%  '

$ws.Range("D23").Value = '8.80'
$ws.Range("E23").Value = '  +7.71%  '

$ws.Range("D24").Value = '5.18'
$ws.Range("E24").Value = '  +5.63%  '

$ws.Range("D25").Value = '6.16'
$ws.Range("E25").Value = '  +12.59%  '

$ws.Range("D26").Value = '12.27'
$ws.Range("E26").Value = '  +4.10%  '

$ws.Range("D27").Value = '86.63'
$ws.Range("E27").Value = '  +5.85%  '

$ws.Range("D28").Value = '3.276.86'
$ws.Range("E28").Value = '  +2.43%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").Value = '  -2.82%  '

$ws.Range("D31").Value = '9.06'
$ws.Range("E31").Value = '  +12.95%  '

$ws.Range("D32").Value = '525.14'
$ws.Range("E32").Value = '  +3.63%  '

$ws.Range("E33").Value = '  -17.10%  '

$ws.Range("E34").Value = '  +4.32%  '

$ws.Range("D35").Value = '7.16'
$ws.Range("E35").Value = '  +6.73%  '

$ws.Range("E36").Value = '  +12.65%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '1.30'
$ws.Range("E37").Value = '  +4.78%  '

$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '23.63'
$ws.Range("E38").Value = '  +6.51%  '

$ws.Range("D40").Value = '0.0906'
$ws.Range("E40").Value = '  +33.29%  '

$ws.Range("D41").Value = '22.28'
$ws.Range("E41").Value = '  +0.21%  '

$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").Value = '0.153'
$ws.Range("E43").Value = '  +16.21%  '

$ws.Range("D44").Value = '0.400'
$ws.Range("E44").Value = '  +11.14%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("E46").Value = '  +6.89%  '

$ws.Range("D47").Value = '148.81'
$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("D48").Value = '44.06'
$ws.Range("E48").Value = '  +1.75%  '

$ws.Range("E49").Value = '  +8.27%  '

$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '4.24'
$ws.Range("E50").Value = '  +8.23%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '167.12'
$ws.Range("E51").Value = '  +5.50%  '
